$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Template")
$ws2 = $wb.Worksheets.Item("Sample Data")

# --- Template sheet: add a blank comment row between the header and the
# sample data row (row 2 becomes blank, old row2/row3 shift down to row3/row4) ---
$ws1.Rows(2).Insert()
$ws1.Rows(2).ClearFormats()
$ws1.Rows(2).ClearContents()

# --- Sample Data sheet: move the header row up from row 3 to row 2 ---
$ws2.Range("A3:P3").Cut($ws2.Range("A2"))
$ws2.Rows(3).Clear()
$ws2.Rows(3).EntireRow.AutoFit()

# --- update selections to match the new layout ---
$ws1.Range("A3:XFD4").Select()
$ws2.Activate()
$ws2.Range("A2:XFD2").Select()
